# Apply the weekly update to the Níspero (Vega Modelo de Temuco) sheet.
# The data rows 2-11 are permuted (row 9 stays fixed) as new weekly
# observations replace the old ones. Only the data columns (D, L, M, N,
# O, P, Q, R, S, T) move; the descriptive columns (A, B, C, E, F, G, H,
# I, J, K) are identical across all rows and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> source row number (values copied from source row)
$rowMap = @{
    2  = 6
    3  = 2
    4  = 8
    5  = 3
    6  = 4
    7  = 11
    8  = 5
    10 = 7
    11 = 10
}

# Columns that carry the data being shuffled between rows.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the current (pre-edit) values of the relevant columns for every
# data row before any writes happen, so that later writes don't clobber
# values still needed as a source for another row.
$snapshot = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 11; $r++) {
        $snapshot["$col$r"] = $ws.Range("$col$r").Value2
    }
}

foreach ($newRow in $rowMap.Keys) {
    $srcRow = $rowMap[$newRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $snapshot["$col$srcRow"]
    }
}
